# Update "想去人数" (F column) counts that changed between scrapes.
# gh-pages data refresh: output generated at 456a3b4

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 246
$ws1.Range("F5").Value = 263
$ws1.Range("F6").Value = 1055
$ws1.Range("F7").Value = 1393
$ws1.Range("F8").Value = 578
$ws1.Range("F14").Value = 397
$ws1.Range("F15").Value = 1269
$ws1.Range("F17").Value = 74
$ws1.Range("F18").Value = 258
$ws1.Range("F19").Value = 5214
$ws1.Range("F20").Value = 629
$ws1.Range("F23").Value = 5462
$ws1.Range("F26").Value = 81
$ws1.Range("F28").Value = 13859
$ws1.Range("F29").Value = 1402
$ws1.Range("F30").Value = 184
$ws1.Range("F31").Value = 81
$ws1.Range("F33").Value = 389
$ws1.Range("F34").Value = 541
$ws1.Range("F35").Value = 4136
$ws1.Range("F36").Value = 81
$ws1.Range("F37").Value = 346

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F4").Value = 40

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 246
$ws4.Range("F5").Value = 263
$ws4.Range("F6").Value = 1055
$ws4.Range("F7").Value = 1393
$ws4.Range("F8").Value = 578
$ws4.Range("F14").Value = 397
$ws4.Range("F15").Value = 1269
$ws4.Range("F17").Value = 74
$ws4.Range("F18").Value = 258
$ws4.Range("F20").Value = 5214
$ws4.Range("F21").Value = 629
$ws4.Range("F25").Value = 40
$ws4.Range("F26").Value = 5462
$ws4.Range("F29").Value = 81
$ws4.Range("F31").Value = 13859
$ws4.Range("F32").Value = 1402
$ws4.Range("F33").Value = 184
$ws4.Range("F34").Value = 81
$ws4.Range("F36").Value = 389
$ws4.Range("F37").Value = 541
$ws4.Range("F38").Value = 4136
$ws4.Range("F39").Value = 81
$ws4.Range("F40").Value = 346
